# ---------------------------------------------------------------------------
# Reorganise the "Formulaire" sheet of the stations workbook:
#  - move chsta_codehydro / chsta_codemeteofrance from the end of the row
#    to just after chsta_suivipluvio (new columns P,Q)
#  - move chsta_infl_ant_type / chsta_infl_nappe from the end of the row
#    to just after chsta_distberge (now columns AF,AG)
#  - update the two named ranges describing the row extent
#  - update the active selection on the sheet
#  - refresh the column widths for the shifted / new columns
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The engine stores ColumnWidth internally with a fixed +5/6 padding versus
# the "raw" width value that ends up in the xlsx xml, so subtract that
# padding to reproduce a precise target width.
$widthPad = 0.8333333333333334

# ----- 1. capture the current contents (value + font size) of columns P:AR -----
$sourceCols = @("P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR")

$values = @{}
$fontSizes = @{}
foreach ($c in $sourceCols) {
    $addr = $c + "1"
    $values[$c] = $ws.Range($addr).Value2
    $fontSizes[$c] = $ws.Range($addr).Font.Size
}

# ----- 2. mapping from the current column letter to its new column letter -----
$colMap = @{
    "P"  = "R"
    "Q"  = "S"
    "R"  = "T"
    "S"  = "U"
    "T"  = "V"
    "U"  = "W"
    "V"  = "X"
    "W"  = "Y"
    "X"  = "Z"
    "Y"  = "AA"
    "Z"  = "AB"
    "AA" = "AC"
    "AB" = "AD"
    "AC" = "AE"
    "AD" = "AH"
    "AE" = "AI"
    "AF" = "AJ"
    "AG" = "AK"
    "AH" = "AL"
    "AI" = "AM"
    "AJ" = "AN"
    "AK" = "AO"
    "AL" = "AP"
    "AM" = "AQ"
    "AN" = "AR"
    "AO" = "P"
    "AP" = "Q"
    "AQ" = "AF"
    "AR" = "AG"
}

# ----- 3. clear the old cells (values + direct formatting) before rewriting -----
foreach ($c in $sourceCols) {
    $addr = $c + "1"
    $cell = $ws.Range($addr)
    $cell.ClearContents()
    $cell.Font.Size = 12
}

# ----- 4. write the values back out at their new location -----
foreach ($c in $sourceCols) {
    $target = $colMap[$c]
    $addr = $target + "1"
    $cell = $ws.Range($addr)
    $cell.Value = $values[$c]
    $cell.Font.Size = $fontSizes[$c]
}

# ----- 5. update the two named ranges that describe the header row extent -----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Formulaire!Modèle_saisie_stations_1") {
        $n.RefersTo = "=Formulaire!`$A`$1:`$AO`$1"
    } elseif ($n.Name -eq "Formulaire!Modèle_saisie_stations") {
        $n.RefersTo = "=Formulaire!`$A`$1:`$AI`$1"
    }
}

# ----- 6. column widths: the final width that should end up on each column  -----
# ----- P..AN, whether it is a shifted column or one of the new columns      -----
$finalWidths = @{
    "P"  = 14.83203125
    "Q"  = 14.83203125
    "R"  = 19.33203125
    "S"  = 18.83203125
    "T"  = 9.1640625
    "U"  = 12.6640625
    "V"  = 12.6640625
    "W"  = 10.83203125
    "X"  = 14.1640625
    "Y"  = 15.1640625
    "Z"  = 16.83203125
    "AA" = 12.6640625
    "AB" = 14.33203125
    "AC" = 13.6640625
    "AD" = 13.83203125
    "AE" = 14.33203125
    "AF" = 14.33203125
    "AG" = 14.33203125
    "AH" = 13.1640625
    "AI" = 13.83203125
    "AJ" = 10.5
    "AK" = 10.5
    "AL" = 10.5
    "AM" = 20.33203125
    "AN" = 14.83203125
}
foreach ($c in $finalWidths.Keys) {
    $addr = $c + "1"
    $ws.Range($addr).EntireColumn.ColumnWidth = [double]$finalWidths[$c] - $widthPad
}

# ----- 7. restore the selection that was active on the sheet -----
$ws.Range("P5").Select()

$wb.Save()
